$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Elements")

# --- Metadata sheet updates ---

# Version 5.0.0 -> 6.0.0
$ws1.Range("B3").Value = "6.0.0"

# Date 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank; set it to "Alvearie Team"
$ws1.Range("B9").Value = "Alvearie Team"

# Remove the duplicate "Contact" / "No display for ContactDetail" row (row 11),
# leaving the single former "Contact" row (row 10) which we turn into the new
# "Jurisdiction" / "United States of America" row.
$ws1.Rows.Item(11).Delete()
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# --- Elements sheet updates ---
# Row 2 (the root Extension element) gets a real Slice Name / Alias instead of
# the generic defaults.
$ws2.Cells.Item(2, 11).Value = "County Code"
$ws2.Cells.Item(2, 12).Value = "The county code for the address"
